# Apply the Alvearie -> LinuxForHealth rebrand edit to the FHIR StructureDefinition workbook.
$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/aca-market-type"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---
$elements = $wb.Worksheets.Item("Elements")
# Constraint(s) for the root "Extension" element row is cleared out.
$elements.Range("AI2").Value = ""
# Fixed Value on Extension.url mirrors the canonical URL (same text as Metadata!B2).
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/aca-market-type"
